$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Notes" column (B) for the five multiclassification-summary rows
# referenced "Table 2" in the supplementary file; the corrected text now
# points to "Table 3". Filenames (column A) are untouched.
$ws.Range("B3").Value = "The file implements the multiclassification machine learning results using the data on Bangladesh at district level for Covid-19 cases. Detailed tables are provided in Section D1 in the supplementary file and summarized in Table 3."
$ws.Range("B7").Value = "The file implements the multiclassification machine learning results using the data on India at district level for Covid-19 cases. Detailed tables are provided in Section D1 in the supplementary file and summarized in Table 3."
$ws.Range("B9").Value = "The file implements the multiclassification machine learning results using the data on India at district level for Covid-19 deaths. Detailed tables are provided in Section D1 in the supplementary file and summarized in Table 3."
$ws.Range("B15").Value = "The file implements the multiclassification machine learning results using the data on Pakistan (Sindh) at district level for Covid-19 cases. Detailed tables are provided in Section D1 in the supplementary file and summarized in Table 3."
$ws.Range("B17").Value = "The file implements the multiclassification machine learning results using the data on Pakistan (Sindh) at district level for Covid-19 deaths. Detailed tables are provided in Section D1 in the supplementary file and summarized in Table 3."

# Move the active selection to B20, matching the author's last-saved cursor
# position for this upload.
$ws.Range("B20").Select()
